$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: new diary entry for 1/30/2020 ---
# Reuse formatting from row 16 (same A/B/C pattern: date, "5:00 -7:50 pm"-style time, "N/A" participants)
$ws.Range("A16:C16").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)

$ws.Range("A18").Value = 43860
$ws.Range("B18").Value = "5:00 -7:50 pm"
$ws.Range("C18").Value = "N/A"
$ws.Range("D18").Value = "UML relationships, overview of visual tools and diagrams."
$ws.Range("E18").Value = "we learned new models to read and understand code"
$ws.Range("F18").Value = "I felt that we went over the UML diagram for a little too long and skipped over sequence graphs and call graphs too quickly. I would have liked to had more time spent on the visual tools because they seemed helpful for reading and understanding larger codebases"
$ws.Range("G18").Value = "Good"

$ws.Rows(18).RowHeight = 136

# --- Row 19: new diary entry for 2/1/2020 ---
# Reuse formatting from row 17 (date style + italic Goal-column style)
$ws.Range("A17").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D19").PasteSpecial(-4122)

$ws.Range("A19").Value = 43862
$ws.Range("A19").NumberFormat = "d-mmm"
$ws.Range("B19").Value = "1:00 - 5:00pm"
$ws.Range("C19").Value = "Chris, Jay, Rafael"
$ws.Range("D19").Value = "Work on open source project: find 2 essential features and create a packet containing everything relevant to those features"
$ws.Range("E19").Value = "We were able to find 2 features and write a packet"
$ws.Range("F19").Value = "I found that going thorugh the codebase to find essential features was much easier now that I was familiar with the program. It was nice knowing where to search and using call graphs to sift through relveant methods and classes"
$ws.Range("G19").Value = "Good"

$ws.Rows(19).RowHeight = 119

# --- Update view / selection state to match where the user ended up editing ---
$ws.Range("A18:XFD19").Select()
